$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 currently holds the "LOT2052 ..." text, row 24 holds "LOT2028 ..." text.
# The edit swaps their order so "LOT2028 ..." comes first (row 23) and
# "LOT2052 ..." comes second (row 24).
$lot2052 = "LOT2052 -  Tecnologia de Bebidas Experimental  (Indicação de Conjunto)`n"
$lot2028 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"

$ws.Range("B23").Value = $lot2028
$ws.Range("C23").Value = $lot2028

$ws.Range("B24").Value = $lot2052
$ws.Range("C24").Value = $lot2052
